# "committing final Gantt chart"
#
# The author cleared out a batch of leftover sample/placeholder task rows
# (their "Assigned To" / "Task" text in columns B & C for rows 13-31,
# excluding the rows that were actually filled in for this project) and
# updated a couple of progress/date values, then left the selection
# sitting on B8 on the ProjectSchedule sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# Rows 13-18, 20, 21 (C only), 22-26, 28-31: clear the leftover
# "Assigned To" (col B) / "Task" (col C) sample text, but keep the cell
# formatting (style) untouched - just drop the contents like the diff
# shows (the <c> elements keep their s="" but lose t="s"/<v>).
$cellsToClear = @(
    "B13", "C13",
    "B14", "C14",
    "B15", "C15",
    "B16", "C16",
    "B17", "C17",
    "B18", "C18",
    "B20", "C20",
    "C21",
    "B22", "C22",
    "B23", "C23",
    "B24", "C24",
    "B25", "C25",
    "B26", "C26",
    "B28", "C28",
    "B29", "C29",
    "B30", "C30",
    "B31", "C31"
)
foreach ($cellAddr in $cellsToClear) {
    $ws.Range($cellAddr).ClearContents()
}

# Progress + dates that changed alongside the cleanup.
$ws.Range("D10").Value = 0.5
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 44620

# Leave the sheet selection where the author left it.
$ws.Activate()
$ws.Range("B8").Select()
